$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: rename "Age" column header to "City"
$ws.Range("C1").Value = "City"

# Row 2: replace numeric Age value with City name (string)
$ws.Range("C2").Value = "Bharatpur"

# Row 3: replace numeric Age value with City name (string)
$ws.Range("C3").Value = "Bangalore"

# Update the selected/active cell to C3
$ws.Range("C3").Select()
